$d = $word.ActiveDocument

# Locate the field token "{m:self.name}" so we don't depend on hard-coded
# character offsets.
$fieldRange = $d.Content.Duplicate
$found = $fieldRange.Find.Execute("{m:self.name}", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$fieldStart = $fieldRange.Start
$fieldEnd = $fieldRange.End

# The run currently holding "{m" must be split right after the opening
# brace, i.e. between "{" and "m".
$openBraceSplit = $d.Range($fieldStart + 1, $fieldStart + 1)

# The run currently holding "name}" must be split right before the closing
# brace, i.e. between "name" and "}".
$closeBraceSplit = $d.Range($fieldEnd - 1, $fieldEnd - 1)

# Word (and this COM host) naturally splits a run in two whenever a
# zero-length structural marker -- such as a bookmark -- is anchored
# strictly inside it. Adding the bookmark and immediately deleting it
# leaves the text/runs themselves untouched but keeps the run boundary
# that was created, which is exactly the "split the run in two" edit
# described by the diff (no extra formatting is introduced).
$bm1 = $d.Bookmarks.Add("m2docSplitPoint1", $openBraceSplit)
$d.Bookmarks("m2docSplitPoint1").Delete()

$bm2 = $d.Bookmarks.Add("m2docSplitPoint2", $closeBraceSplit)
$d.Bookmarks("m2docSplitPoint2").Delete()
